# Add "Aaron Schomer, " before "Bandon Cho, ..." in the team-roster paragraph.
$d = $word.ActiveDocument

# Locate the paragraph that currently starts the roster line.
$rng = $d.Content
$rng.Find.Execute("Bandon Cho, Jaceguai de Magalhaes, Jack Nguyen, ") | Out-Null
$para = $rng.Paragraphs(1)
$paraRng = $para.Range

# Pull the paragraph's own opening <w:p .../> tag (with its w14:paraId / rsid
# attributes) out of its OpenXML so the replacement keeps the same identity.
$owx = $paraRng.WordOpenXML
$pOpenTag = "<w:p>"
if ($owx -match "<w:p( [^>]*)?>") {
    $pOpenTag = $matches[0]
}

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# New runs: "Aaron " / "Schomer" (flagged spellStart/spellEnd) / ", "
$newRuns = ""
$newRuns += "<w:r><w:t xml:space=`"preserve`">Aaron </w:t></w:r>"
$newRuns += "<w:proofErr w:type=`"spellStart`"/>"
$newRuns += "<w:r><w:t>Schomer</w:t></w:r>"
$newRuns += "<w:proofErr w:type=`"spellEnd`"/>"
$newRuns += "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>"

# Keep the rest of the paragraph exactly as-is (existing runs + proofErr marks).
$existingInner = "<w:r><w:t xml:space=`"preserve`">Bandon Cho, Jaceguai de Magalhaes, Jack Nguyen, </w:t></w:r>"
$existingInner += "<w:proofErr w:type=`"spellStart`"/>"
$existingInner += "<w:r><w:t>Newyork</w:t></w:r>"
$existingInner += "<w:proofErr w:type=`"spellEnd`"/>"
$existingInner += "<w:r><w:t xml:space=`"preserve`"> Her</w:t></w:r>"

$packageXml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
    "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=`"$wNs`"><w:body>$pOpenTag$newRuns$existingInner</w:p></w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

$paraRng.InsertXML($packageXml)
